# Update namespace and create course model
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The old 3-row "Gos.Rab" news layout is replaced by a 2-row, 3-column
# course-card layout (title / description / icon path). Drop row 3 first.
$ws.Range("A3:B3").ClearContents() | Out-Null

# Row 1 - "Yandex Lyceum" course card
$ws.Range("A1").Value = "Яндекс лицей"
$ws.Range("B1").Value = "Комампания яндекс в очередной раз проводит делает набор в ряды своих юных разработчиков.
                        Бесплатное обучение проходит в течении 4 полугодии, в течении которых ученики получают хорошие
                        знания языка програмирования Python."
$ws.Range("C1").Value = "../static/icon/Yandex_lyceum.jpg"

# Row 2 - "GitHub" course card
$ws.Range("A2").Value = "GitHub"
$ws.Range("B2").Value = "GitHub — это крупнейший веб-сервис для хостинга IT-проектов и их совместной разработки.
                    Веб-сервис основан на системе контроля версий Git и разработан на Ruby on Rails и Erlang компанией
                    GitHub, Inc (ранее Logical Awesome). Сервис бесплатен для проектов с открытым исходным кодом и (с
                    2019 года) небольших частных проектов, предоставляя им все возможности (включая SSL), а для крупных
                    корпоративных проектов предлагаются различные платные тарифные планы."
$ws.Range("C2").Value = "../static/icon/git.png"

# Descriptions wrap, matching the new cell style (cellXfs index 1)
$ws.Range("B1:B2").WrapText = $true

# Column layout: narrower title column, wider description/icon columns
$ws.Columns.Item(1).ColumnWidth = 14.83
$ws.Range("B1:C2").EntireColumn.ColumnWidth = 29.83

# Tall rows to fit the wrapped course descriptions
$ws.Rows.Item(1).RowHeight = 150
$ws.Rows.Item(2).RowHeight = 315

# Page setup for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where the author last clicked
$ws.Range("E2").Select() | Out-Null

